$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '246.26'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.358'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.05813'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '3.377'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '6.475'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.8099'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.9209'
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.01068'
$ws.Range('E10').Value = '9OneONEBestin24h'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.1402'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07401'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.03190'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.03034'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.09373'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.851'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.001552'
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.04712'
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.006053'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.004691'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.00008798'
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.3184'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.03839'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.006389'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1066'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.003099'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.009063'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005257'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.7098'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0002000'
